# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 case tracker sorted by total cases
# (column B) descending. This refresh (11:42 -> 12:12) updates the
# numbers for several countries; since the table stays sorted by total
# cases, a handful of neighbouring rows swap places as their counts
# cross over each other.
#
# Each row below is written with its final country name (column A) and
# its final Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes values (columns B-H), in the
# row position it occupies after the refreshed data is re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row 6   'Estados Unidos' 68581 370  428 67117 1455 9  1036
Set-Row 8   'Alemania'       39457 2134 3547 35688 23  16 222
Set-Row 11  'Suiza'          11125 228  131 10829 141 12 165
Set-Row 16  'Austria'        6001  413  112 5847  28  11 42
Set-Row 21  'Brasil'         2563  9    6   2497  18  1  60
Set-Row 22  'Suecia'         2554  28   16  2474  176 2  64
Set-Row 34  'Polonia'        1085  34   7   1064  3   0  14

# Rumania overtakes Finlandia in total cases -> rows 36/37 swap.
Set-Row 36  'Rumania'        1029  123  94  918   29  0  17
Set-Row 37  'Finlandia'      915   35   10  900   22  2  5

Set-Row 46  'India'          694   37   45  635   0   2  14

# Libano overtakes Irak in total cases -> rows 62/63 swap.
Set-Row 62  'Libano'         368   35   20  342   3   0  6
Set-Row 63  'Irak'           346   0    103 214   0   0  29

Set-Row 86  'Moldavia'       149   0    2   146   28  0  1
Set-Row 98  'Sri Lanka'      102   0    7   95    3   0  0

# Bolivia overtakes Puerto Rico in total cases -> rows 123/124 swap.
Set-Row 123 'Bolivia'        39    7    0   39    0   0  0
Set-Row 124 'Puerto Rico'    39    0    1   36    0   0  2

# Haiti overtakes Surinam in total cases -> rows 153/154 swap.
Set-Row 153 'Haiti'          8     0    0   8     0   0  0
Set-Row 154 'Surinam'        8     0    0   8     0   0  0

# Gabon/Namibia/Bermudas reshuffle around Niger (row 158 stays Niger).
Set-Row 157 'Gabon'          7     1    0   6     0   0  1
Set-Row 159 'Namibia'        7     0    2   5     0   0  0
Set-Row 160 'Bermudas'       7     0    2   5     0   0  0

# Footer timestamp refresh.
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 12:12"

Write-Host "Update complete"
